$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.04349144772309899
$ws.Range("D2").Value = 0.01883115997351759
$ws.Range("E2").Value = 0.07309914496272718
$ws.Range("F2").Value = 3.621130701086102
$ws.Range("G2").Value = 0.002595240643110432
$ws.Range("I2").Value = 0.848719136981444
$ws.Range("J2").Value = 0.1532231040505181
$ws.Range("K2").Value = 2.372091779716072
$ws.Range("M2").Value = 0.6127212123230592
$ws.Range("N2").Value = 2.671174786595792

# Row 3
$ws.Range("B3").Value = 0.03785547096567399
$ws.Range("D3").Value = 0.01872192813357287
$ws.Range("E3").Value = 0.07300949670209178
$ws.Range("F3").Value = 3.60153998104623
$ws.Range("G3").Value = 0.002600610571142101
$ws.Range("I3").Value = 0.8520888305691372
$ws.Range("J3").Value = 0.153165599077326
$ws.Range("K3").Value = 2.270071448469082
$ws.Range("M3").Value = 0.5943551562412708
$ws.Range("N3").Value = 2.689019173921089

# Row 4
$ws.Range("B4").Value = 0.0343869474186107
$ws.Range("D4").Value = 0.01865813169924913
$ws.Range("E4").Value = 0.07298061153166557
$ws.Range("F4").Value = 3.591331531968194
$ws.Range("G4").Value = 0.002604081324794145
$ws.Range("I4").Value = 0.8544607495672274
$ws.Range("J4").Value = 0.153192703655435
$ws.Range("K4").Value = 2.208896038737947
$ws.Range("M4").Value = 0.5834405332797203
$ws.Range("N4").Value = 2.700725547982636

# Row 5
$ws.Range("B5").Value = 0.03297160571764834
$ws.Range("D5").Value = 0.01863296198187747
$ws.Range("E5").Value = 0.07297542873228302
$ws.Range("F5").Value = 3.587628259199036
$ws.Range("G5").Value = 0.002605539486467722
$ws.Range("I5").Value = 0.8555034485468411
$ws.Range("J5").Value = 0.1532194367017361
$ws.Range("K5").Value = 2.184333476566422
$ws.Range("M5").Value = 0.5790835811395141
$ws.Range("N5").Value = 2.705684334839049

# Row 6
$ws.Range("B6").Value = 0.03273647880862285
$ws.Range("D6").Value = 0.01862883272085725
$ws.Range("E6").Value = 0.07297496634447143
$ws.Range("F6").Value = 3.587040891961394
$ws.Range("G6").Value = 0.002605784263063016
$ws.Range("I6").Value = 0.8556811841809271
$ws.Range("J6").Value = 0.1532248230245621
$ws.Range("K6").Value = 2.180276996651457
$ws.Range("M6").Value = 0.578365592664035
$ws.Range("N6").Value = 2.706519107058526

# Row 7
$ws.Range("B7").Value = 0.03436786711371553
$ws.Range("D7").Value = 0.01865778889277259
$ws.Range("E7").Value = 0.07298051494555935
$ws.Range("F7").Value = 3.591279740283809
$ws.Range("G7").Value = 0.002604100812369118
$ws.Range("I7").Value = 0.8544745036202279
$ws.Range("J7").Value = 0.1531930006757172
$ws.Range("K7").Value = 2.208563295913933
$ws.Range("M7").Value = 0.5833814063793312
$ws.Range("N7").Value = 2.7007916616176

# Row 8
$ws.Range("B8").Value = 0.04154989827220845
$ws.Range("D8").Value = 0.01879282052084186
$ws.Range("E8").Value = 0.0730628107778557
$ws.Range("F8").Value = 3.613997398512367
$ws.Range("G8").Value = 0.002597056254684299
$ws.Range("I8").Value = 0.8498181082881615
$ws.Range("J8").Value = 0.1531903181336318
$ws.Range("K8").Value = 2.336610002350426
$ws.Range("M8").Value = 0.606313268910732
$ws.Range("N8").Value = 2.67717177349688

# Row 9
$ws.Range("B9").Value = 0.05556549742021843
$ws.Range("D9").Value = 0.01908339626983668
$ws.Range("E9").Value = 0.07343139099512719
$ws.Range("F9").Value = 3.673041885308749
$ws.Range("G9").Value = 0.002584612475586283
$ws.Range("I9").Value = 0.8430929071706998
$ws.Range("J9").Value = 0.1536807548546193
$ws.Range("K9").Value = 2.599426075254598
$ws.Range("M9").Value = 0.6541705389353467
$ws.Range("N9").Value = 2.636812233843045

# Row 10
$ws.Range("B10").Value = 0.06581567295077662
$ws.Range("D10").Value = 0.01931239821938213
$ws.Range("E10").Value = 0.07382813753561024
$ws.Range("F10").Value = 3.725340739528349
$ws.Range("G10").Value = 0.002576295998853466
$ws.Range("I10").Value = 0.8396226238529678
$ws.Range("J10").Value = 0.1543442017487635
$ws.Range("K10").Value = 2.79980758306516
$ws.Range("M10").Value = 0.6911160955265387
$ws.Range("N10").Value = 2.610805361767632

# Row 11
$ws.Range("B11").Value = 0.07046743732998095
$ws.Range("D11").Value = 0.01941991146091127
$ws.Range("E11").Value = 0.07403591229672735
$ws.Range("F11").Value = 3.75108824080246
$ws.Range("G11").Value = 0.00257268992537973
$ws.Range("I11").Value = 0.8383641754610451
$ws.Range("J11").Value = 0.1547120699650364
$ws.Range("K11").Value = 2.892583251735346
$ws.Range("M11").Value = 0.7083167799823116
$ws.Range("N11").Value = 2.599768506601762

# Row 12
$ws.Range("B12").Value = 0.07222723259461361
$ws.Range("D12").Value = 0.01946110091047437
$ws.Range("E12").Value = 0.07411850875259063
$ws.Range("F12").Value = 3.761120865394417
$ws.Range("G12").Value = 0.002571349714487033
$ws.Range("I12").Value = 0.8379337503797615
$ws.Range("J12").Value = 0.1548608871491695
$ws.Range("K12").Value = 2.927950410464121
$ws.Range("M12").Value = 0.7148872562117532
$ws.Range("N12").Value = 2.595703551228411

# Row 13
$ws.Range("B13").Value = 0.07184830809843845
$ws.Range("D13").Value = 0.01945220888662647
$ws.Range("E13").Value = 0.07410054605808369
$ws.Range("F13").Value = 3.758947570849188
$ws.Range("G13").Value = 0.002571637228702454
$ws.Range("I13").Value = 0.838024397580817
$ws.Range("J13").Value = 0.1548284133662321
$ws.Range("K13").Value = 2.920322976258149
$ws.Range("M13").Value = 0.7134696488691787
$ws.Range("N13").Value = 2.596573916277492

# Row 14
$ws.Range("B14").Value = 0.07061225210956934
$ws.Range("D14").Value = 0.01942329061098036
$ws.Range("E14").Value = 0.07404262910606363
$ws.Range("F14").Value = 3.751907959716277
$ws.Range("G14").Value = 0.002572579158431898
$ws.Range("I14").Value = 0.8383278393274765
$ws.Range("J14").Value = 0.1547241224992746
$ws.Range("K14").Value = 2.895488207783899
$ws.Range("M14").Value = 0.7088561939065983
$ws.Range("N14").Value = 2.599431784460307

# Row 15
$ws.Range("B15").Value = 0.06985490332174038
$ws.Range("D15").Value = 0.01940563928596362
$ws.Range("E15").Value = 0.07400766313263674
$ws.Range("F15").Value = 3.747632836054692
$ws.Range("G15").Value = 0.002573159412678105
$ws.Range("I15").Value = 0.8385197149328221
$ws.Range("J15").Value = 0.1546614806852844
$ws.Range("K15").Value = 2.88030685986962
$ws.Range("M15").Value = 0.7060377451496578
$ws.Range("N15").Value = 2.601197227038583

# Row 16
$ws.Range("B16").Value = 0.06551143471712351
$ws.Range("D16").Value = 0.01930543884944314
$ws.Range("E16").Value = 0.0738151074766602
$ws.Range("F16").Value = 3.723697556778092
$ws.Range("G16").Value = 0.002576535216621858
$ws.Range("I16").Value = 0.839711316265948
$ws.Range("J16").Value = 0.154321491207412
$ws.Range("K16").Value = 2.793777258395608
$ws.Range("M16").Value = 0.6899999432702586
$ws.Range("N16").Value = 2.611542646246946

# Row 17
$ws.Range("B17").Value = 0.06284392241911974
$ws.Range("D17").Value = 0.01924482167270902
$ws.Range("E17").Value = 0.07370396537479351
$ws.Range("F17").Value = 3.70951603309922
$ws.Range("G17").Value = 0.002578651430758638
$ws.Range("I17").Value = 0.8405243908719982
$ws.Range("J17").Value = 0.154129848516618
$ws.Range("K17").Value = 2.741110704199002
$ws.Range("M17").Value = 0.6802623850939824
$ws.Range("N17").Value = 2.618092758979984

# Row 18
$ws.Range("B18").Value = 0.06130860372205404
$ws.Range("D18").Value = 0.01921027076808102
$ws.Range("E18").Value = 0.07364260876507878
$ws.Range("F18").Value = 3.701543272760063
$ws.Range("G18").Value = 0.002579885302268388
$ws.Range("I18").Value = 0.8410221839841299
$ws.Range("J18").Value = 0.1540258380035269
$ws.Range("K18").Value = 2.710970765381944
$ws.Range("M18").Value = 0.6746986627162812
$ws.Range("N18").Value = 2.62193490772421

# Row 19
$ws.Range("B19").Value = 0.06078859725914754
$ws.Range("D19").Value = 0.0191986265556352
$ws.Range("E19").Value = 0.07362227601414162
$ws.Range("F19").Value = 3.698875410398017
$ws.Range("G19").Value = 0.002580305939365616
$ws.Range("I19").Value = 0.841195901398649
$ws.Range("J19").Value = 0.1539916892494801
$ws.Range("K19").Value = 2.70079202929918
$ws.Range("M19").Value = 0.6728212386508829
$ws.Range("N19").Value = 2.623248613092144

# Row 20
$ws.Range("B20").Value = 0.06312799185540996
$ws.Range("D20").Value = 0.01925124195305727
$ws.Range("E20").Value = 0.07371553077820892
$ws.Range("F20").Value = 3.711006620170025
$ws.Range("G20").Value = 0.002578424430793378
$ws.Range("I20").Value = 0.8404347183713625
$ws.Range("J20").Value = 0.1541496056444629
$ws.Range("K20").Value = 2.746701355005257
$ws.Range("M20").Value = 0.6812951267143319
$ws.Range("N20").Value = 2.617387755173567

# Row 21
$ws.Range("B21").Value = 0.07097535963261237
$ws.Range("D21").Value = 0.01943177170489285
$ws.Range("E21").Value = 0.07405953448413172
$ws.Range("F21").Value = 3.753967984134732
$ws.Range("G21").Value = 0.002572301804346529
$ws.Range("I21").Value = 0.8382374587669474
$ws.Range("J21").Value = 0.1547544969494439
$ws.Range("K21").Value = 2.902776387285371
$ws.Range("M21").Value = 0.7102097294241076
$ws.Range("N21").Value = 2.59858925028783

# Row 22
$ws.Range("B22").Value = 0.07609394136984804
$ws.Range("D22").Value = 0.01955253435504289
$ws.Range("E22").Value = 0.07430718657328583
$ws.Range("F22").Value = 3.783693560618076
$ws.Range("G22").Value = 0.002568447897539012
$ws.Range("I22").Value = 0.8370702940706209
$ws.Range("J22").Value = 0.1552052870248204
$ws.Range("K22").Value = 3.006151275634295
$ws.Range("M22").Value = 0.729439188467623
$ws.Range("N22").Value = 2.586970664445246

# Row 23
$ws.Range("B23").Value = 0.07336302798947258
$ws.Range("D23").Value = 0.0194878281806794
$ws.Range("E23").Value = 0.07417292400892173
$ws.Range("F23").Value = 3.767677258887346
$ws.Range("G23").Value = 0.002570491342618653
$ws.Range("I23").Value = 0.8376686040505987
$ws.Range("J23").Value = 0.1549596127068398
$ws.Range("K23").Value = 2.950852077985815
$ws.Range("M23").Value = 0.7191455767461292
$ws.Range("N23").Value = 2.593110556483467

# Row 24
$ws.Range("B24").Value = 0.06299956944269525
$ws.Range("D24").Value = 0.01924833841369278
$ws.Range("E24").Value = 0.07371029414577812
$ws.Range("F24").Value = 3.710332163986664
$ws.Range("G24").Value = 0.00257852700375139
$ws.Range("I24").Value = 0.8404751647947108
$ws.Range("J24").Value = 0.1541406542378994
$ws.Range("K24").Value = 2.744173389698233
$ws.Range("M24").Value = 0.6808281166983221
$ws.Range("N24").Value = 2.617706249444041

# Row 25
$ws.Range("B25").Value = 0.05178178611335227
$ws.Range("D25").Value = 0.01900204966551655
$ws.Range("E25").Value = 0.0733095367225296
$ws.Range("F25").Value = 3.65550856867165
$ws.Range("G25").Value = 0.002587833109829111
$ws.Range("I25").Value = 0.8446542706177667
$ws.Range("J25").Value = 0.1534949167924751
$ws.Range("K25").Value = 2.52705787989396
$ws.Range("M25").Value = 0.6409122264747467
$ws.Range("N25").Value = 2.647091338075796
